$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("L1").Value = "From Currency"
$ws.Range("M1").Value = "To Currency"
$ws.Range("N1").Value = "Exchange Rate "
$ws.Range("O1").Value = "As Of"

# --- Row 2 ---
$ws.Range("L2").Value = "USD"
$ws.Range("M2").Value = "INR"
$ws.Range("N2").Value = 80
$ws.Range("O2").Value2 = 44743

# --- Row 3 ---
$ws.Range("L3").Value = "USD"
$ws.Range("M3").Value = "INR"
$ws.Range("N3").Value = 81
$ws.Range("O3").Value2 = 44774

# --- Row 4 ---
$ws.Range("L4").Value = "USD"
$ws.Range("M4").Value = "INR"
$ws.Range("N4").Value = 81
$ws.Range("O4").Value2 = 44866

# Copy the date number format (style) from column G onto the new "As Of" column
# so it matches the existing date formatting (reuses the same style index).
$ws.Range("G2").Copy()
$ws.Range("O2:O4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the sheet view: scroll so column C is the left-most visible
# column, and move the active selection to the newly added O4 cell ---
$excel.ActiveWindow.ScrollColumn = 3
[void]$ws.Range("O4").Select()
